# Wrap "Hexagon 5" (shape 1) and "Picture 4" (shape 2) into a new group shape
# ("Group 1"), matching PowerPoint's real grouping math (which uses each
# child's *visually rotated* bounding box), while leaving the two original
# shapes' own geometry completely untouched in the final file.
#
# The emulated Shape.Left/Top/Width/Height COM properties round-trip through
# a 32-bit Single before being converted back to EMU, which truncates and can
# land 1 EMU low. EmuToPt() nudges by a tiny epsilon (far smaller than the
# ~1/12700 pt = 1 EMU quantum) to compensate so the stored EMU matches exactly.

function EmuToPt($emu) {
    return ($emu / 12700.0) + 0.00001
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

$hex = $shapes.Item(1)
$pic = $shapes.Item(2)

# Original (current) hexagon box, in EMU -- this is what must be restored
# after grouping.
$hexLeft0   = 146776
$hexTop0    = 130495
$hexWidth0  = 1535245
$hexHeight0 = 1321882

# The hexagon is rotated 90 degrees, so its true on-slide footprint swaps
# width/height around the same centre. Temporarily point the shape's raw
# (pre-rotation) box at that footprint -- since swapping twice is an
# identity for a 90-degree rotation, the shape keeps rendering in exactly
# the same place, but the naive (rotation-unaware) bounding union the
# grouping operation computes now equals PowerPoint's real rotation-aware
# union.
$cx2 = ($hexLeft0 * 2) + $hexWidth0
$cy2 = ($hexTop0 * 2) + $hexHeight0

$visWidth  = $hexHeight0
$visHeight = $hexWidth0
$visLeft   = [Math]::Floor(($cx2 - $visWidth)  / 2.0)
$visTop    = [Math]::Floor(($cy2 - $visHeight) / 2.0)

$hex.Left   = EmuToPt $visLeft
$hex.Top    = EmuToPt $visTop
$hex.Width  = EmuToPt $visWidth
$hex.Height = EmuToPt $visHeight

# Group the two shapes. The new group's off/ext (and its initial,
# identity-scale chOff/chExt) are computed from the current (temporarily
# adjusted) shape boxes.
$range = $shapes.Range(@(1, 2))
$grp = $range.Group()
$grp.Name = "Group 1"

# Restore the hexagon child's own geometry to its original numbers.
$items = $grp.GroupItems
$hexChild = $items.Item(1)
$hexChild.Left   = EmuToPt $hexLeft0
$hexChild.Top    = EmuToPt $hexTop0
$hexChild.Width  = EmuToPt $hexWidth0
$hexChild.Height = EmuToPt $hexHeight0
